$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("产品记录")

# ---------------------------------------------------------------------------
# Header row (row 1) text updates
# ---------------------------------------------------------------------------
$ws.Range("H1").Value = "支付金额"
$ws.Range("I1").Value = "每年交费次数"
$ws.Range("J1").Value = "交费期间（年）"
$ws.Range("K1").Value = "保险期间（年）"
$ws.Range("O1").Value = "返还比例（%）"

# ---------------------------------------------------------------------------
# Existing row edits
# ---------------------------------------------------------------------------
# Row 2 - 一次支付 / 到期返还带利率
$ws.Range("L2").Value = 1100
$ws.Range("O2").Value = 1.1

# Row 3 - 分期付款 / 每月支付，为期一年
$ws.Range("I3").Value = 12

# Row 4 - 年金保险 / 前五年每年付款，后二十年每月返还
$ws.Range("I4").Value = 1
$ws.Range("O4").Value = 1.1

# Row 5 - 2013年7月15日开始扣款...
$ws.Range("I5").Value = 12

# ---------------------------------------------------------------------------
# New row 6 - 都来保险 / 分期付满期返
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "都来保险"
$ws.Range("D6").Value = "分期付满期返"
$ws.Range("E6").Value = 2019
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 28
$ws.Range("H6").Value = 720
$ws.Range("I6").Value = 12
$ws.Range("J6").Value = 10
$ws.Range("K6").Value = 20
$ws.Range("L6").Value = 475200
$ws.Range("O6").Value = 1.1

# ---------------------------------------------------------------------------
# New row 7 - 意外险 / 消费型
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "意外险"
$ws.Range("D7").Value = "消费型"
$ws.Range("E7").Value = 2020
$ws.Range("F7").Value = 6
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 397
$ws.Range("I7").Value = 12
$ws.Range("J7").Value = 10

# ---------------------------------------------------------------------------
# Formatting: A6 / A7 get a non-bold red font (flagging the two new entries);
# D7 drops the wrap-text formatting used by the other rows in column D.
# ---------------------------------------------------------------------------
$ws.Range("A6").Font.Bold = $false
$ws.Range("A6").Font.Color = 255
$ws.Range("A7").Font.Bold = $false
$ws.Range("A7").Font.Color = 255
$ws.Range("D7").Font.Bold = $false
$ws.Range("D7").WrapText = $false

# ---------------------------------------------------------------------------
# View: freeze pane / selection moved to reflect the extra columns and rows
# ---------------------------------------------------------------------------
$ws.Range("M17").Select()

Write-Output "done"
